$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Sent By" column header
$ws.Range("H1").Value = "Sent By"

# Update Date Time column (D) - both rows now share the same new date
$ws.Range("D2").Value = "15 February 2025, 12:00 AM"
$ws.Range("D3").Value = "15 February 2025, 12:00 AM"

# Update Email Subject column (E)
$ws.Range("E2").Value = "Welcome to Our Institute"
$ws.Range("E3").Value = "Exam Schedule Announcement"

# Update Status column (G)
$ws.Range("G2").Value = "Pending"
$ws.Range("G3").Value = "Pending"

# Populate new "Sent By" column (H)
$ws.Range("H2").Value = "John Smith"
$ws.Range("H3").Value = "John Smith"
